$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 4519
$ws.Range("B3").Value = 4890
$ws.Range("B4").Value = 3517
$ws.Range("B5").Value = 1257

$ws.Range("B19").Select()
